$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.946.06'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '1.817.59'
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4644'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3649'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07221'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.22%  '
$ws.Range('E10').Value = '  -3.11%  '
$ws.Range('E11').Value = '  -3.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07564'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.99%  '
$ws.Range('D13').Value = '1.778.82'
$ws.Range('E13').Value = '  -5.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.318'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.79'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.471'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008613'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.40%  '
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.37%  '
$ws.Range('D21').Value = '26.679.71'
$ws.Range('E21').Value = '  -3.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.133'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.43%  '
$ws.Range('E23').Value = '  -1.72%  '
$ws.Range('D24').Value = '2.012.46'
$ws.Range('E24').Value = '  -4.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.87'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.851'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.066'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.104'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.20'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08864'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.974'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.413'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.127'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7142'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.074'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05246'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.412'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.40%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01921'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.916'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.127'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5142'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1622'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.149'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4804'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.71'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06254'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.612'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.91'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.98%  '
